$p = $ppt.ActivePresentation

# Slide 2: "Git" section title
$s2 = $p.Slides.Add(2, 2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Git"

# Slide 3: "Create local repository"
$s3 = $p.Slides.Add(3, 2)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Create local repository"
$tf3 = $s3.Shapes.Item(2).TextFrame
$tr3 = $tf3.TextRange
$r3 = $tr3
$tr3.Text = "mkdir"
$r3 = $r3.InsertAfter(" project(create project folder)")
$r3 = $r3.InsertAfter("`rcd project(move into project folder)")
$r3 = $r3.InsertAfter("`rgit")
$r3 = $r3.InsertAfter(" ")
$r3 = $r3.InsertAfter("init")
$r3 = $r3.InsertAfter("(initial local ")
$r3 = $r3.InsertAfter("git")
$r3 = $r3.InsertAfter(")")
$r3 = $r3.InsertAfter("`rtouch test.txt(create test.txt)")
$r3 = $r3.InsertAfter("`rg")
$r3 = $r3.InsertAfter("it")
$r3 = $r3.InsertAfter(" status(check  untrack file)")
$r3 = $r3.InsertAfter("`rg")
$r3 = $r3.InsertAfter("it")
$r3 = $r3.InsertAfter(" add .(track file)")
$r3 = $r3.InsertAfter("`rgit")
$r3 = $r3.InsertAfter(" commit  –m “commit description”")
$r3 = $r3.InsertAfter("`rg")
$r3 = $r3.InsertAfter("it")
$r3 = $r3.InsertAfter(" log(check commit history)")
for ($i = 1; $i -le $tr3.Paragraphs().Count; $i++) {
  $tr3.Paragraphs($i, 1).ParagraphFormat.Bullet.Visible = 0
}
$tf3.AutoSize = 2

# Slide 4: empty separator slide
$s4 = $p.Slides.Add(4, 2)

# Slide 5: "Github" section title
$s5 = $p.Slides.Add(5, 2)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Github"

# Slide 6: git clone / push / pull
$s6 = $p.Slides.Add(6, 2)
$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange
$tr6.Text = "git"
$r6 = $tr6
$r6 = $r6.InsertAfter(" clone (clone repository to local) ")
$r6 = $r6.InsertAfter("`rgit")
$r6 = $r6.InsertAfter(" push(push local repository to ")
$r6 = $r6.InsertAfter("github")
$r6 = $r6.InsertAfter(")")
$r6 = $r6.InsertAfter("`rgit")
$r6 = $r6.InsertAfter(" pull (pull repository from ")
$r6 = $r6.InsertAfter("github")
$r6 = $r6.InsertAfter(")")

# Slide 7: empty closing slide
$s7 = $p.Slides.Add(7, 2)

